$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds a date serial number that needs to move from 45206 (2023-10-07)
# to 45208 (2023-10-09) for every data row (rows 2 through 501).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 501 }

$ws.Range("C2:C$lastRow").Value = 45208
